$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.09743279987498
$ws.Range("C2").Value = 4.383610461228727
$ws.Range("D2").Value = 10.55600213707838
$ws.Range("F2").Value = 33.86179476280677
$ws.Range("G2").Value = 3.657505479244528
$ws.Range("J2").Value = 11.38325827074975
$ws.Range("K2").Value = 10.4597486080256
$ws.Range("N2").Value = 19.14434137877175
$ws.Range("O2").Value = 24.97602651856182

$ws.Range("B3").Value = 10.83781369857153
$ws.Range("C3").Value = 4.168377122335583
$ws.Range("D3").Value = 10.4589168366868
$ws.Range("F3").Value = 33.88778119270932
$ws.Range("G3").Value = 3.659363388589991
$ws.Range("J3").Value = 11.36493578756304
$ws.Range("K3").Value = 10.28343398044409
$ws.Range("N3").Value = 19.2055266257816
$ws.Range("O3").Value = 25.04599247561742

$ws.Range("B4").Value = 10.67736933921746
$ws.Range("C4").Value = 4.029562607868116
$ws.Range("D4").Value = 10.40108082809133
$ws.Range("F4").Value = 33.91182969454054
$ws.Range("G4").Value = 3.660564825462544
$ws.Range("J4").Value = 11.35595289567657
$ws.Range("K4").Value = 10.17543649536271
$ws.Range("N4").Value = 19.24483376828429
$ws.Range("O4").Value = 25.09402202378385

$ws.Range("B5").Value = 10.61182511364603
$ws.Range("C5").Value = 3.97135813303171
$ws.Range("D5").Value = 10.37798206355747
$ws.Range("F5").Value = 33.92366252279353
$ws.Range("G5").Value = 3.661069725444963
$ws.Range("J5").Value = 11.35286503494215
$ws.Range("K5").Value = 10.13154850109054
$ws.Range("N5").Value = 19.26129040188989
$ws.Range("O5").Value = 25.11486695088306

$ws.Range("B6").Value = 10.60093473193928
$ws.Range("C6").Value = 3.961595749928219
$ws.Range("D6").Value = 10.37417558274318
$ws.Range("F6").Value = 33.92575005507791
$ws.Range("G6").Value = 3.661154489480253
$ws.Range("J6").Value = 11.35238695856799
$ws.Range("K6").Value = 10.12426995028228
$ws.Range("N6").Value = 19.2640495460941
$ws.Range("O6").Value = 25.11840501707322

$ws.Range("B7").Value = 10.67648590877792
$ws.Range("C7").Value = 4.02878421191021
$ws.Range("D7").Value = 10.400767377043
$ws.Range("F7").Value = 33.91198104871093
$ws.Range("G7").Value = 3.660571572688746
$ws.Range("J7").Value = 11.35590892954376
$ws.Range("K7").Value = 10.17484403827408
$ws.Range("N7").Value = 19.24505393035091
$ws.Range("O7").Value = 25.0942979968087

$ws.Range("B8").Value = 11.00818960608062
$ws.Range("C8").Value = 4.310801263874475
$ws.Range("D8").Value = 10.5221731991363
$ws.Range("F8").Value = 33.86907423756725
$ws.Range("G8").Value = 3.658133522627886
$ws.Range("J8").Value = 11.37647183549302
$ws.Range("K8").Value = 10.39893490986299
$ws.Range("N8").Value = 19.16507789472773
$ws.Range("O8").Value = 24.99909703739859

$ws.Range("B9").Value = 11.64628544038969
$ws.Range("C9").Value = 4.809634038481338
$ws.Range("D9").Value = 10.77325527133033
$ws.Range("F9").Value = 33.84920357238433
$ws.Range("G9").Value = 3.653831737042217
$ws.Range("J9").Value = 11.43465059022827
$ws.Range("K9").Value = 10.8380044240158
$ws.Range("N9").Value = 19.02198280643023
$ws.Range("O9").Value = 24.8527426574427

$ws.Range("B10").Value = 12.10233017445237
$ws.Range("C10").Value = 5.141658488313791
$ws.Range("D10").Value = 10.96416751669497
$ws.Range("F10").Value = 33.87380825567511
$ws.Range("G10").Value = 3.650960295688212
$ws.Range("J10").Value = 11.48807322981178
$ws.Range("K10").Value = 11.15718046164337
$ws.Range("N10").Value = 18.92513778969998
$ws.Range("O10").Value = 24.76994217850792

$ws.Range("B11").Value = 12.30602499108403
$ws.Range("C11").Value = 5.28501047667866
$ws.Range("D11").Value = 11.05209243960369
$ws.Range("F11").Value = 33.89349831362332
$ws.Range("G11").Value = 3.649716121523666
$ws.Range("J11").Value = 11.51464239227978
$ws.Range("K11").Value = 11.30100045910224
$ws.Range("N11").Value = 18.88286114099818
$ws.Range("O11").Value = 24.73766965083033

$ws.Range("B12").Value = 12.38254218857024
$ws.Range("C12").Value = 5.338175660657723
$ws.Range("D12").Value = 11.08551481229582
$ws.Range("F12").Value = 33.90217308500073
$ws.Range("G12").Value = 3.649253859754977
$ws.Range("J12").Value = 11.52502423085542
$ws.Range("K12").Value = 11.35521294736712
$ws.Range("N12").Value = 18.86710645249206
$ws.Range("O12").Value = 24.72622642526742

$ws.Range("B13").Value = 12.3660914676535
$ws.Range("C13").Value = 5.326775545669305
$ws.Range("D13").Value = 11.07831149660125
$ws.Range("F13").Value = 33.90025068587668
$ws.Range("G13").Value = 3.649353021855562
$ws.Range("J13").Value = 11.52277414307047
$ws.Range("K13").Value = 11.34354918517062
$ws.Range("N13").Value = 18.87048820395545
$ws.Range("O13").Value = 24.72865631027379

$ws.Range("B14").Value = 12.31233290302199
$ws.Range("C14").Value = 5.289406904576551
$ws.Range("D14").Value = 11.05483973041859
$ws.Range("F14").Value = 33.89418759315906
$ws.Range("G14").Value = 3.649677913252862
$ws.Range("J14").Value = 11.51549011749878
$ws.Range("K14").Value = 11.30546582943346
$ws.Range("N14").Value = 18.88155989837767
$ws.Range("O14").Value = 24.73671261276928

$ws.Range("B15").Value = 12.27932162411685
$ws.Range("C15").Value = 5.266371438141641
$ws.Range("D15").Value = 11.04047832244331
$ws.Range("F15").Value = 33.89063235085919
$ws.Range("G15").Value = 3.649878073750783
$ws.Range("J15").Value = 11.51107003799781
$ws.Range("K15").Value = 11.28210476158315
$ws.Range("N15").Value = 18.88837474666185
$ws.Range("O15").Value = 24.74174866554339

$ws.Range("B16").Value = 12.08893585325568
$ws.Range("C16").Value = 5.132134170565143
$ws.Range("D16").Value = 10.95844072248537
$ws.Range("F16").Value = 33.87269219163237
$ws.Range("G16").Value = 3.651042850409773
$ws.Range("J16").Value = 11.48638204479395
$ws.Range("K16").Value = 11.14774938042509
$ws.Range("N16").Value = 18.92793635557342
$ws.Range("O16").Value = 24.77215998132052

$ws.Range("B17").Value = 11.97112041654443
$ws.Range("C17").Value = 5.047804306859859
$ws.Range("D17").Value = 10.90836955719051
$ws.Range("F17").Value = 33.86386082021403
$ws.Range("G17").Value = 3.651773266510614
$ws.Range("J17").Value = 11.47181373294579
$ws.Range("K17").Value = 11.06493732373863
$ws.Range("N17").Value = 18.95266081177597
$ws.Range("O17").Value = 24.79219934436598

$ws.Range("B18").Value = 11.90300619040758
$ws.Range("C18").Value = 4.998577012790501
$ws.Range("D18").Value = 10.87967308560461
$ws.Range("F18").Value = 33.85958131943828
$ws.Range("G18").Value = 3.652199226419993
$ws.Range("J18").Value = 11.46364814475557
$ws.Range("K18").Value = 11.01717954432035
$ws.Range("N18").Value = 18.96704914990252
$ws.Range("O18").Value = 24.80423303633749

$ws.Range("B19").Value = 11.87988621644477
$ws.Range("C19").Value = 4.981785788131318
$ws.Range("D19").Value = 10.86997551785119
$ws.Range("F19").Value = 33.85826985030867
$ws.Range("G19").Value = 3.652344454166178
$ws.Range("J19").Value = 11.4609202783323
$ws.Range("K19").Value = 11.00098946133645
$ws.Range("N19").Value = 18.97194959655209
$ws.Range("O19").Value = 24.80839454963209

$ws.Range("B20").Value = 11.98369885780358
$ws.Range("C20").Value = 5.056856299417576
$ws.Range("D20").Value = 10.91368924597978
$ws.Range("F20").Value = 33.86471815106387
$ws.Range("G20").Value = 3.651694907969612
$ws.Range("J20").Value = 11.4733424725862
$ws.Range("K20").Value = 11.07376627756334
$ws.Range("N20").Value = 18.95001152589951
$ws.Range("O20").Value = 24.79001357182888

$ws.Range("B21").Value = 12.32814044016612
$ws.Range("C21").Value = 5.300413446614882
$ws.Range("D21").Value = 11.06173072940378
$ws.Range("F21").Value = 33.8959354306033
$ws.Range("G21").Value = 3.64958224412249
$ws.Range("J21").Value = 11.51762095376011
$ws.Range("K21").Value = 11.31665898978394
$ws.Range("N21").Value = 18.87830097418378
$ws.Range("O21").Value = 24.7343251599264

$ws.Range("B22").Value = 12.54962156629734
$ws.Range("C22").Value = 5.453064776797067
$ws.Range("D22").Value = 11.15921243304262
$ws.Range("F22").Value = 33.92343819649909
$ws.Range("G22").Value = 3.648253235349651
$ws.Range("J22").Value = 11.54842587896479
$ws.Range("K22").Value = 11.47392983325954
$ws.Range("N22").Value = 18.83291720710886
$ws.Range("O22").Value = 24.70246330534703

$ws.Range("B23").Value = 12.43176900707567
$ws.Range("C23").Value = 5.372192911632123
$ws.Range("D23").Value = 11.10712719464136
$ws.Range("F23").Value = 33.90811113882782
$ws.Range("G23").Value = 3.648957832654204
$ws.Range("J23").Value = 11.53181579716438
$ws.Range("K23").Value = 11.39014258876482
$ws.Range("N23").Value = 18.85700405791359
$ws.Range("O23").Value = 24.71905310660787

$ws.Range("B24").Value = 11.97801332397768
$ws.Range("C24").Value = 5.052766212848005
$ws.Range("D24").Value = 10.91128393187895
$ws.Range("F24").Value = 33.86432806651152
$ws.Range("G24").Value = 3.651730315067572
$ws.Range("J24").Value = 11.47265067509678
$ws.Range("K24").Value = 11.06977516465811
$ws.Range("N24").Value = 18.95120872618631
$ws.Range("O24").Value = 24.7910001623965

$ws.Range("B25").Value = 11.47554724073412
$ws.Range("C25").Value = 4.680655104892271
$ws.Range("D25").Value = 10.70409007049546
$ws.Range("F25").Value = 33.8476921591525
$ws.Range("G25").Value = 3.654944498217679
$ws.Range("J25").Value = 11.4170197160862
$ws.Range("K25").Value = 10.71959875170883
$ws.Range("N25").Value = 19.05923223362324
$ws.Range("O25").Value = 24.88800257353786

